# Youth voter turnout.xlsx - "Add files via upload" edit
#
# 1. Edit the "Young Adults (TUFTS)" sheet:
#    - remove the stray "democrat"/"Republican" cells in C1:D1
#    - remove the stray "Percent increase: " cell in row 7
#    - insert a new bold title row at the top
# 2. Reorder the sheet tabs so "Young Adults (TUFTS)" comes before
#    "College students - TUFTS"
# 3. Make "Young Adults (TUFTS)" the active/selected tab, with D2 selected

$wb = $excel.ActiveWorkbook

$wsYoung = $wb.Worksheets.Item("Young Adults (TUFTS)")

# Drop the orphan "democrat"/"Republican" header cells that lived outside
# the real table, and the leftover "Percent increase: " label in row 7.
$wsYoung.Range("C1:D1").ClearContents()
$wsYoung.Rows.Item(7).Delete()

# Insert a bold title row above the existing header row.
$wsYoung.Rows.Item(1).Insert()
$wsYoung.Range("A1").Value = "Youth Voter Turnout in Recent Midterm Elections "
$wsYoung.Range("A1").Font.Bold = $true

# Move "Young Adults (TUFTS)" so it sits before "College students - TUFTS".
$wsCollege = $wb.Worksheets.Item("College students - TUFTS")
$wsYoung.Move($wsCollege)

# Re-fetch after the move and make it the active sheet / selection.
$wsYoung = $wb.Worksheets.Item("Young Adults (TUFTS)")
$wsYoung.Activate()
$wsYoung.Range("D2").Select()
